$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (dates) must stay as literal text, not be auto-converted to Excel date serials,
# so force Text format on the destination range before writing the date strings.
$ws.Range("A634:A654").NumberFormat = "@"

$rows = @(
    @{ "A" = "2024-08-28"; "C" = 1890.599975585938; "D" = 711.7999877929688; "E" = 78.95999908447266; "F" = 299.9500122070312; "G" = 1323.25; "H" = 29275.72984313965; "I" = 0; "J" = 353.4546403129226 },
    @{ "A" = "2024-08-29"; "C" = 1893.25; "D" = 700.7999877929688; "E" = 77; "F" = 296.2000122070312; "G" = 1298.650024414062; "H" = 28928.15008544922; "I" = -0.01187262485180639; "J" = 349.2582059663571 },
    @{ "A" = "2024-08-30"; "C" = 1899.349975585938; "D" = 696.0999755859375; "E" = 75.83999633789062; "F" = 299.2999877929688; "G" = 1302.900024414062; "H" = 28869.96923828125; "I" = -0.002011219071945895; "J" = 348.5557712014839 },
    @{ "A" = "2024-09-02"; "C" = 1885.400024414062; "D" = 683.5999755859375; "E" = 73.80999755859375; "F" = 296.8999938964844; "G" = 1303.849975585938; "H" = 28487.57955932617; "I" = -0.01324524026329871; "J" = 343.9390662667609 },
    @{ "A" = "2024-09-03"; "C" = 1901.949951171875; "D" = 689.4000244140625; "E" = 74.47000122070312; "F" = 297.1499938964844; "G" = 1320.25; "H" = 28715.86001586914; "I" = 0.008013332830455756; "J" = 346.6951644781526 },
    @{ "A" = "2024-09-04"; "C" = 1901.300048828125; "D" = 688.9500122070312; "E" = 74.16000366210938; "F" = 298.9500122070312; "G" = 1327.099975585938; "H" = 28721.03076171875; "I" = 0.0001800658537390795; "J" = 346.7575924389315 },
    @{ "A" = "2024-09-05"; "C" = 1879.449951171875; "D" = 687.5; "E" = 76; "F" = 290.6000061035156; "G" = 1312.349975585938; "H" = 28602.74978637695; "I" = -0.004118270556621158; "J" = 345.3295508557054 },
    @{ "A" = "2024-09-06"; "C" = 1872.349975585938; "D" = 673.5499877929688; "E" = 74.72000122070312; "F" = 283.6000061035156; "G" = 1289.699951171875; "H" = 28191.60983276367; "I" = -0.01437414083205038; "J" = 340.3657352582368 },
    @{ "A" = "2024-09-09"; "C" = 1892.400024414062; "D" = 664.1500244140625; "E" = 74.33999633789062; "F" = 281.5499877929688; "G" = 1237.150024414062; "H" = 28036.46997070312; "I" = -0.005503050836077006; "J" = 338.492685314252 },
    @{ "A" = "2024-09-10"; "C" = 1922.449951171875; "D" = 664.5999755859375; "E" = 78.05000305175781; "F" = 285.75; "G" = 1250.300048828125; "H" = 28561.49984741211; "I" = 0.01872667555001102; "J" = 344.831528008184 },
    @{ "A" = "2024-09-11"; "C" = 1957.599975585938; "D" = 689.75; "E" = 81.94999694824219; "F" = 288.0499877929688; "G" = 1237.699951171875; "H" = 29297.64938354492; "I" = 0.0257741904334731; "J" = 353.7192814785324 },
    @{ "A" = "2024-09-12"; "C" = 1996.400024414062; "D" = 729.1500244140625; "E" = 81.69999694824219; "F" = 291.7000122070312; "G" = 1237.300048828125; "H" = 29912.80038452148; "I" = 0.02099659917843318; "J" = 361.1461834534205 },
    @{ "A" = "2024-09-13"; "C" = 1988.050048828125; "D" = 713.7000122070312; "E" = 83.11000061035156; "F" = 289.9500122070312; "G" = 1241.5; "H" = 29812.18057250977; "I" = -0.003363771051799782; "J" = 359.9313703760519 },
    @{ "A" = "2024-09-16"; "C" = 1989.900024414062; "D" = 714.2000122070312; "E" = 84.69999694824219; "F" = 290.3999938964844; "G" = 1226.599975585938; "H" = 29926.49987792969; "I" = 0.003834650911961043; "J" = 361.3115815337079 },
    @{ "A" = "2024-09-17"; "C" = 2006.550048828125; "D" = 731.0999755859375; "E" = 82; "F" = 284.2999877929688; "G" = 1193.800048828125; "H" = 29823.24993896484; "I" = -0.003450117433913109; "J" = 360.0650141471837 },
    @{ "A" = "2024-09-18"; "C" = 1987.800048828125; "D" = 743.25; "E" = 80.81999969482422; "F" = 282.8500061035156; "G" = 1166.400024414062; "H" = 29685.31034851074; "I" = -0.004625236710834788; "J" = 358.3996282254629 },
    @{ "A" = "2024-09-19"; "C" = 1998.599975585938; "D" = 735.9500122070312; "E" = 80.97000122070312; "F" = 272.7000122070312; "G" = 1121.300048828125; "H" = 29455.8603515625; "I" = -0.0077294120982553; "J" = 355.6294098030468 },
    @{ "A" = "2024-09-20"; "C" = 2048.10009765625; "D" = 746.5; "E" = 83.44999694824219; "F" = 277.3500061035156; "G" = 1149.400024414062; "H" = 30118.95037841797; "I" = 0.02251131078642199; "J" = 363.6350939719151 },
    @{ "A" = "2024-09-23"; "C" = 2082.39990234375; "D" = 773.9500122070312; "E" = 82.88999938964844; "F" = 286.2999877929688; "G" = 1162.75; "H" = 30664.31942749023; "I" = 0.01810717313253569; "J" = 370.2194975755305 },
    @{ "A" = "2024-09-24"; "C" = 2068.14990234375; "D" = 781.8499755859375; "E" = 83.79000091552734; "F" = 291.7999877929688; "G" = 1141.199951171875; "H" = 30770.6690826416; "I" = 0.003468188994144963; "J" = 371.5034887624398 },
    @{ "A" = "2024-09-25"; "C" = 2061.60009765625; "D" = 775.8499755859375; "E" = 82.95999908447266; "F" = 289.8500061035156; "G" = 1118.449951171875; "H" = 30542.33015441895; "I" = -0.007420668286718119; "J" = 368.7466846049752 },
)

$r = 634
foreach ($row in $rows) {
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
    $r++
}

# Drop the Text number format now that the literal date strings are safely
# stored (prevents Excel from re-parsing "2024-08-28" as a date serial),
# so the cells end up plain/general-formatted text like the rest of column A.
$ws.Range("A634:A654").ClearFormats()

Write-Output "Added rows 634-654"